$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H17").Value = 1367.6666
$ws.Range("J17").Value = 1367.6666
$ws.Range("L17").Value = 4102.9998
$ws.Range("N17").Value = -4438.9998
$ws.Range("H40").Value = 4023.7334
$ws.Range("I40").Value = 8586.666999999999
$ws.Range("J40").Value = 2883
$ws.Range("K40").Value = 8586.666999999999
$ws.Range("L40").Value = 2883
$ws.Range("M40").Value = -8411.666999999999
$ws.Range("N40").Value = -3233
$ws.Range("H107").Value = 713.86365
$ws.Range("I107").Value = 735.75
$ws.Range("J107").Value = 495
$ws.Range("K107").Value = 735.75
$ws.Range("L107").Value = 495
$ws.Range("M107").Value = 1184.25
$ws.Range("N107").Value = -4335
$ws.Range("H129").Value = 2420.6875
$ws.Range("I129").Value = 12999.625
$ws.Range("J129").Value = 909.4107
$ws.Range("K129").Value = 38998.875
$ws.Range("L129").Value = 2728.2321
$ws.Range("M129").Value = -33998.875
$ws.Range("N129").Value = -12728.2321
$ws.Range("H137").Value = 1305.3334
$ws.Range("I137").Value = 1191
$ws.Range("K137").Value = 3573
$ws.Range("M137").Value = -1023
$ws.Range("H138").Value = 3790.3115
$ws.Range("I138").Value = 2706.2144
$ws.Range("J138").Value = 4113.234
$ws.Range("K138").Value = 8118.6432
$ws.Range("L138").Value = 12339.702
$ws.Range("M138").Value = -2978.6432
$ws.Range("N138").Value = -22619.702

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 31554.871
$ws.Range("I32").Value = 11418.646
$ws.Range("K32").Value = 11418.646
$ws.Range("M32").Value = -11131.646
$ws.Range("H45").Value = 53708.156
$ws.Range("I45").Value = 83925.914
$ws.Range("K45").Value = 83925.914
$ws.Range("M45").Value = -83548.914
$ws.Range("H74").Value = 834.4314000000001
$ws.Range("I74").Value = 761.5217
$ws.Range("K74").Value = 761.5217
$ws.Range("M74").Value = 112.4783
$ws.Range("H77").Value = 834.4314000000001
$ws.Range("I77").Value = 761.5217
$ws.Range("K77").Value = 3807.6085
$ws.Range("M77").Value = 560.3914999999997
$ws.Range("H132").Value = 14410.702
$ws.Range("I132").Value = 16044.049
$ws.Range("K132").Value = 48132.147
$ws.Range("M132").Value = -45602.147

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H86").Value = 93508.336
$ws.Range("I86").Value = 101863.63
$ws.Range("J86").Value = 1600
$ws.Range("K86").Value = 101863.63
$ws.Range("L86").Value = 1600
$ws.Range("M86").Value = -100740.63
$ws.Range("N86").Value = -3846
$ws.Range("H89").Value = 93508.336
$ws.Range("I89").Value = 101863.63
$ws.Range("J89").Value = 1600
$ws.Range("K89").Value = 509318.15
$ws.Range("L89").Value = 8000
$ws.Range("M89").Value = -503702.15
$ws.Range("N89").Value = -19232
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 334.54544
$ws.Range("I94").Value = 306.19232
$ws.Range("J94").Value = 439.85715
$ws.Range("K94").Value = 306.19232
$ws.Range("L94").Value = 439.85715
$ws.Range("M94").Value = 144.80768
$ws.Range("N94").Value = -1341.85715
$ws.Range("H134").Value = 3271.4897
$ws.Range("I134").Value = 3712.1
$ws.Range("J134").Value = 2575.7896
$ws.Range("K134").Value = 11136.3
$ws.Range("L134").Value = 7727.3688
$ws.Range("M134").Value = -8601.299999999999
$ws.Range("N134").Value = -12797.3688

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H25").Value = 19866.666
$ws.Range("J25").Value = 19866.666
$ws.Range("L25").Value = 19866.666
$ws.Range("N25").Value = -20214.666
$ws.Range("H31").Value = 50323.5
$ws.Range("I31").Value = 1388.5333
$ws.Range("J31").Value = 99258.47
$ws.Range("K31").Value = 1388.5333
$ws.Range("L31").Value = 99258.47
$ws.Range("M31").Value = -1093.5333
$ws.Range("N31").Value = -99848.47
$ws.Range("H34").Value = 50323.5
$ws.Range("I34").Value = 1388.5333
$ws.Range("J34").Value = 99258.47
$ws.Range("K34").Value = 1388.5333
$ws.Range("L34").Value = 99258.47
$ws.Range("M34").Value = -1186.5333
$ws.Range("N34").Value = -99662.47
$ws.Range("H62").Value = 2388.24
$ws.Range("I62").Value = 2240
$ws.Range("J62").Value = 2610.6
$ws.Range("K62").Value = 2240
$ws.Range("L62").Value = 2610.6
$ws.Range("M62").Value = -1616
$ws.Range("N62").Value = -3858.6
$ws.Range("H65").Value = 2388.24
$ws.Range("I65").Value = 2240
$ws.Range("J65").Value = 2610.6
$ws.Range("K65").Value = 11200
$ws.Range("L65").Value = 13053
$ws.Range("M65").Value = -8080
$ws.Range("N65").Value = -19293
$ws.Range("H74").Value = 38544.332
$ws.Range("J74").Value = 38544.332
$ws.Range("L74").Value = 38544.332
$ws.Range("N74").Value = -40292.332
$ws.Range("H77").Value = 38544.332
$ws.Range("J77").Value = 38544.332
$ws.Range("L77").Value = 115632.996
$ws.Range("N77").Value = -124368.996
$ws.Range("H122").Value = 1385.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1385.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 4156.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9056.5
$ws.Range("H132").Value = 4147.1
$ws.Range("I132").Value = 4784.6665
$ws.Range("K132").Value = 14353.9995
$ws.Range("M132").Value = -11823.9995

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H12").Value = 50.2
$ws.Range("J12").Value = 56.058823
$ws.Range("L12").Value = 168.176469
$ws.Range("N12").Value = -514.176469
$ws.Range("H14").Value = 652.55554
$ws.Range("I14").Value = 652.55554
$ws.Range("K14").Value = 1957.66662
$ws.Range("M14").Value = -1784.66662
$ws.Range("H39").Value = 4380
$ws.Range("J39").Value = 5350
$ws.Range("L39").Value = 16050
$ws.Range("N39").Value = -16638
$ws.Range("H137").Value = 4046684.5
$ws.Range("I137").Value = 73440.64
$ws.Range("J137").Value = 9103541
$ws.Range("K137").Value = 220321.92
$ws.Range("L137").Value = 27310623
$ws.Range("M137").Value = -215221.92
$ws.Range("N137").Value = -27320823

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H43").Value = 2030.5151
$ws.Range("I43").Value = 998.0769
$ws.Range("J43").Value = 5865.2856
$ws.Range("K43").Value = 998.0769
$ws.Range("L43").Value = 5865.2856
$ws.Range("M43").Value = -847.0769
$ws.Range("N43").Value = -6167.2856
$ws.Range("H122").Value = 2445.8667
$ws.Range("I122").Value = 2492.923
$ws.Range("J122").Value = 2140
$ws.Range("K122").Value = 7478.768999999999
$ws.Range("L122").Value = 6420
$ws.Range("M122").Value = -5028.768999999999
$ws.Range("N122").Value = -11320
$ws.Range("H126").Value = 2549.8235
$ws.Range("I126").Value = 2187.2307
$ws.Range("J126").Value = 3728.25
$ws.Range("K126").Value = 6561.6921
$ws.Range("L126").Value = 11184.75
$ws.Range("M126").Value = -4091.6921
$ws.Range("N126").Value = -16124.75
$ws.Range("H132").Value = 2246.4717
$ws.Range("I132").Value = 1802.2258
$ws.Range("J132").Value = 2872.4546
$ws.Range("K132").Value = 5406.6774
$ws.Range("L132").Value = 8617.363799999999
$ws.Range("M132").Value = -2876.6774
$ws.Range("N132").Value = -13677.3638

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H51").Value = 19799.834
$ws.Range("J51").Value = 19799.834
$ws.Range("L51").Value = 19799.834
$ws.Range("N51").Value = -20755.834
$ws.Range("H68").Value = 4198.5454
$ws.Range("I68").Value = 2457.7144
$ws.Range("J68").Value = 7245
$ws.Range("K68").Value = 2457.7144
$ws.Range("L68").Value = 7245
$ws.Range("M68").Value = -1708.7144
$ws.Range("N68").Value = -8743
$ws.Range("H71").Value = 4198.5454
$ws.Range("I71").Value = 2457.7144
$ws.Range("J71").Value = 7245
$ws.Range("K71").Value = 12288.572
$ws.Range("L71").Value = 36225
$ws.Range("M71").Value = -8544.572
$ws.Range("N71").Value = -43713
$ws.Range("H122").Value = 3527.6667
$ws.Range("I122").Value = 3329.8
$ws.Range("J122").Value = 3775
$ws.Range("K122").Value = 9989.400000000001
$ws.Range("L122").Value = 11325
$ws.Range("M122").Value = -7539.400000000001
$ws.Range("N122").Value = -16225

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H62").Value = 10991367
$ws.Range("I62").Value = 76923070
$ws.Range("J62").Value = 2750
$ws.Range("K62").Value = 76923070
$ws.Range("L62").Value = 2750
$ws.Range("M62").Value = -76922446
$ws.Range("N62").Value = -3998
$ws.Range("H65").Value = 10991367
$ws.Range("I65").Value = 76923070
$ws.Range("J65").Value = 2750
$ws.Range("K65").Value = 384615350
$ws.Range("L65").Value = 13750
$ws.Range("M65").Value = -384612230
$ws.Range("N65").Value = -19990
$ws.Range("H81").Value = 253909.88
$ws.Range("I81").Value = 1000000
$ws.Range("J81").Value = 147325.58
$ws.Range("K81").Value = 2000000
$ws.Range("L81").Value = 294651.16
$ws.Range("M81").Value = -1998939
$ws.Range("N81").Value = -296773.16
$ws.Range("H84").Value = 253909.88
$ws.Range("I84").Value = 1000000
$ws.Range("J84").Value = 147325.58
$ws.Range("K84").Value = 10000000
$ws.Range("L84").Value = 1473255.8
$ws.Range("M84").Value = -9994696
$ws.Range("N84").Value = -1483863.8
$ws.Range("H100").Value = 100639.8
$ws.Range("I100").Value = 143228.42
$ws.Range("J100").Value = 1266.3334
$ws.Range("K100").Value = 286456.84
$ws.Range("L100").Value = 2532.6668
$ws.Range("M100").Value = -285915.84
$ws.Range("N100").Value = -3614.6668
$ws.Range("H122").Value = 1497
$ws.Range("I122").Value = 1549.3158
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4647.9474
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2197.9474
$ws.Range("N122").Value = -7900
